# "start process of adding games"
#
# The workbook's four generic sheets are being repurposed/labelled for the
# different phases of a game (Ranges / Kickoff & PAT / Punt / Field Goal),
# and the author leaves off having the "Field Goal" sheet active (with a
# specific cell selected) ready to keep working on it.

$wb = $excel.ActiveWorkbook

# Give the sheets their real names (was Sheet1..Sheet4).
$wb.Worksheets.Item(1).Name = "Ranges"
$wb.Worksheets.Item(2).Name = "Kickoff and PAT"
$wb.Worksheets.Item(3).Name = "Punt"
$wb.Worksheets.Item(4).Name = "Field Goal"

# Work shifts to the "Field Goal" sheet: make it the active tab and leave
# the selection on H14 (was K33).
$wsFieldGoal = $wb.Worksheets.Item(4)
$wsFieldGoal.Activate()
$wsFieldGoal.Range("H14").Select()
